$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 17: fill in the remaining columns for "33. Search in Rotated Sorted Array" ---
$ws.Range("C17").Value = "Binary Search"
$ws.Range("D17").Value = "Visualize graphically the pivot point, look for the most ideal conditions to narrow the search, and reject the rest."
$ws.Range("E17").Value = "https://leetcode.com/problems/search-in-rotated-sorted-array/solutions/14436/revised-binary-search/ "
$ws.Hyperlinks.Add($ws.Range("E17"), "https://leetcode.com/problems/search-in-rotated-sorted-array/solutions/14436/revised-binary-search/ ", "", "", "https://leetcode.com/problems/search-in-rotated-sorted-array/solutions/14436/revised-binary-search/ ") | Out-Null
$ws.Range("E17").Style = "Hyperlink"

# --- Row 19: new entry, "287. Find the Duplicate Number" ---
$ws.Range("A19").Value = "287. Find the Duplicate Number"
$ws.Range("B19").Value = "Medium"
$ws.Range("B19").Style = $ws.Range("B18").Style
$ws.Range("C19").Value = "Linked List"
$ws.Range("D19").Value = "Linked List Cycle and Floyd's algorithm for cycle start searching (2 phases). Consider the elements as index pointers. "
$ws.Range("E19").Value = "https://leetcode.com/problems/find-the-duplicate-number/solutions/1892921/9-approaches-count-hash-in-place-marked-sort-binary-search-bit-mask-fast-slow-pointers/ "
$ws.Hyperlinks.Add($ws.Range("E19"), "https://leetcode.com/problems/find-the-duplicate-number/solutions/1892921/9-approaches-count-hash-in-place-marked-sort-binary-search-bit-mask-fast-slow-pointers/ ", "", "", "https://leetcode.com/problems/find-the-duplicate-number/solutions/1892921/9-approaches-count-hash-in-place-marked-sort-binary-search-bit-mask-fast-slow-pointers/ ") | Out-Null
$ws.Range("E19").Style = "Hyperlink"

# --- Resize the table to include the new row ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E19"))

# --- Update the view: scroll right one column and move the selection ---
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("E27").Select() | Out-Null
